# Update COVID demographic data across the four sheets
# (Age Group, Gender, Race, Ethnicity) with refreshed counts/rates.

$wb = $excel.ActiveWorkbook

# --- Sheet: Age Group ---
$ws = $wb.Worksheets.Item("Age Group")
$ws.Range("B2").Value = 436
$ws.Range("D2").Value = 2.44

$ws.Range("B3").Value = 2150
$ws.Range("D3").Value = 12.05
$ws.Range("E3").Value = 0.2

$ws.Range("B4").Value = 2614
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 14.66
$ws.Range("E4").Value = 0.89

$ws.Range("B5").Value = 3093
$ws.Range("D5").Value = 17.34
$ws.Range("E5").Value = 2.09

$ws.Range("B6").Value = 3310
$ws.Range("C6").Value = 61
$ws.Range("D6").Value = 18.56
$ws.Range("E6").Value = 6.06

$ws.Range("B7").Value = 2624
$ws.Range("C7").Value = 177
$ws.Range("D7").Value = 14.71
$ws.Range("E7").Value = 17.58

$ws.Range("B8").Value = 1691
$ws.Range("C8").Value = 274
$ws.Range("D8").Value = 9.48
$ws.Range("E8").Value = 27.21

$ws.Range("B9").Value = 1899
$ws.Range("C9").Value = 462
$ws.Range("D9").Value = 10.65
$ws.Range("E9").Value = 45.88

$ws.Range("B10").Value = 18
$ws.Range("D10").Value = 0.1

# --- Sheet: Gender ---
$ws = $wb.Worksheets.Item("Gender")
$ws.Range("B2").Value = 9404
$ws.Range("C2").Value = 449
$ws.Range("D2").Value = 52.73
$ws.Range("E2").Value = 44.59

$ws.Range("B3").Value = 8079
$ws.Range("C3").Value = 536
$ws.Range("D3").Value = 45.3
$ws.Range("E3").Value = 53.23

$ws.Range("B4").Value = 352
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 1.97
$ws.Range("E4").Value = 2.18

# --- Sheet: Race ---
$ws = $wb.Worksheets.Item("Race")
$ws.Range("B2").Value = 278
$ws.Range("D2").Value = 1.56
$ws.Range("E2").Value = 0.4

$ws.Range("B3").Value = 2767
$ws.Range("C3").Value = 171
$ws.Range("D3").Value = 15.51
$ws.Range("E3").Value = 16.98

$ws.Range("B4").Value = 2657
$ws.Range("C4").Value = 112
$ws.Range("D4").Value = 14.9
$ws.Range("E4").Value = 11.12

$ws.Range("B5").Value = 3789
$ws.Range("C5").Value = 68
$ws.Range("D5").Value = 21.24
$ws.Range("E5").Value = 6.75

$ws.Range("B6").Value = 8344
$ws.Range("C6").Value = 652
$ws.Range("D6").Value = 46.78
$ws.Range("E6").Value = 64.75

# --- Sheet: Ethnicity ---
$ws = $wb.Worksheets.Item("Ethnicity")
$ws.Range("B2").Value = 1642
$ws.Range("D2").Value = 9.210000000000001
$ws.Range("E2").Value = 1.39

$ws.Range("B3").Value = 6801
$ws.Range("C3").Value = 557
$ws.Range("D3").Value = 38.13
$ws.Range("E3").Value = 55.31

$ws.Range("B4").Value = 9392
$ws.Range("C4").Value = 436
$ws.Range("D4").Value = 52.66
$ws.Range("E4").Value = 43.3

$wb.Save()
